$p = $ppt.ActivePresentation

# --- Change 1: Slide 10, "Content Placeholder 2" -----------------------
# The runs "Zeslabení " and "napětí" are merged back into a single run
# "Zeslabení napětí" (same rPr: lang="cs-CZ" sz="3000" dirty="0" smtClean="0").
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$full10 = $tr10.Text
$old1 = "Zeslabení napětí"
$idx1 = $full10.IndexOf("Zeslabení")
$rng1 = $tr10.Characters($idx1 + 1, $old1.Length)
$rng1.Text = "Zeslabení napětí"

# --- Change 2: Slide 9, Title -------------------------------------------
# "mikroprocesoru ESP" is replaced by "a jejich řešení" so the title reads
# "Problémy a jejich řešení".
$s9 = $p.Slides.Item(9)
$sh9Title = $s9.Shapes.Item(1)
$tr9Title = $sh9Title.TextFrame.TextRange
$full9Title = $tr9Title.Text
$old2 = "mikroprocesoru ESP"
$idx2 = $full9Title.IndexOf($old2)
$rng2 = $tr9Title.Characters($idx2 + 1, $old2.Length)
$rng2.Text = "a jejich řešení"

# --- Change 3: Slide 9, "Content Placeholder 2" --------------------------
# The typo "ryhlost" is fixed to "rychlost" (run ends up split into
# "Seriová " + "rychlost").
$sh9Body = $s9.Shapes.Item(2)
$tr9Body = $sh9Body.TextFrame.TextRange
$full9Body = $tr9Body.Text
$old3 = "ryhlost"
$idx3 = $full9Body.IndexOf($old3)
$rng3 = $tr9Body.Characters($idx3 + 1, $old3.Length)
$rng3.Text = "rychlost"
